$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.158605
$ws.Range("H2").Value = 0.475815
$ws.Range("I2").Value = 0.005248964619762229
$ws.Range("J2").Value = 0.005248964619762229
$ws.Range("M2").Value = 1.275643
$ws.Range("N2").Value = 3.826929
$ws.Range("O2").Value = 0.008652234199457187
$ws.Range("P2").Value = 0.008652234199457187
$ws.Range("Q2").Value = 0.202323358015
$ws.Range("R2").Value = 1.820910222135
$ws.Range("S2").Value = 0.00004541527119484754
$ws.Range("T2").Value = 0.00004541527119484755
$ws.Range("G3").Value = 0.158605
$ws.Range("H3").Value = 0.475815
$ws.Range("I3").Value = 0.005248964619762229
$ws.Range("J3").Value = 0.005248964619762229
$ws.Range("O3").Value = 0.6643867693241158
$ws.Range("P3").Value = 0.6643867693241158
$ws.Range("Q3").Value = 15.53598285617667
$ws.Range("R3").Value = 139.82384570559
$ws.Range("S3").Value = 0.003487342646020413
$ws.Range("T3").Value = 0.003487342646020413
$ws.Range("G4").Value = 0.158605
$ws.Range("H4").Value = 0.475815
$ws.Range("I4").Value = 0.005248964619762229
$ws.Range("J4").Value = 0.005248964619762229
$ws.Range("M4").Value = 48.20552666666666
$ws.Range("N4").Value = 144.61658
$ws.Range("O4").Value = 0.326960996476427
$ws.Range("P4").Value = 0.326960996476427
$ws.Range("Q4").Value = 7.645637556966666
$ws.Range("R4").Value = 68.8107380127
$ws.Range("S4").Value = 0.001716206702546968
$ws.Range("T4").Value = 0.001716206702546968
$ws.Range("I5").Value = 0.9259451473351599
$ws.Range("J5").Value = 0.92594514733516
$ws.Range("M5").Value = 1.275643
$ws.Range("N5").Value = 3.826929
$ws.Range("O5").Value = 0.008652234199457187
$ws.Range("P5").Value = 0.008652234199457187
$ws.Range("Q5").Value = 35.69091146875167
$ws.Range("R5").Value = 321.218203218765
$ws.Range("S5").Value = 0.008011494270594693
$ws.Range("T5").Value = 0.008011494270594695
$ws.Range("I6").Value = 0.9259451473351599
$ws.Range("J6").Value = 0.92594514733516
$ws.Range("O6").Value = 0.6643867693241158
$ws.Range("P6").Value = 0.6643867693241158
$ws.Range("S6").Value = 0.6151857050093492
$ws.Range("T6").Value = 0.6151857050093493
$ws.Range("I7").Value = 0.9259451473351599
$ws.Range("J7").Value = 0.92594514733516
$ws.Range("M7").Value = 48.20552666666666
$ws.Range("N7").Value = 144.61658
$ws.Range("O7").Value = 0.326960996476427
$ws.Range("P7").Value = 0.326960996476427
$ws.Range("Q7").Value = 1348.730941622811
$ws.Range("R7").Value = 12138.5784746053
$ws.Range("S7").Value = 0.3027479480552158
$ws.Range("T7").Value = 0.3027479480552159
$ws.Range("G8").Value = 2.079068666666667
$ws.Range("H8").Value = 6.237206
$ws.Range("I8").Value = 0.0688058880450778
$ws.Range("J8").Value = 0.06880588804507781
$ws.Range("M8").Value = 1.275643
$ws.Range("N8").Value = 3.826929
$ws.Range("O8").Value = 0.008652234199457187
$ws.Range("P8").Value = 0.008652234199457187
$ws.Range("Q8").Value = 2.652149391152667
$ws.Range("R8").Value = 23.869344520374
$ws.Range("S8").Value = 0.0005953246576676446
$ws.Range("T8").Value = 0.0005953246576676447
$ws.Range("G9").Value = 2.079068666666667
$ws.Range("H9").Value = 6.237206
$ws.Range("I9").Value = 0.0688058880450778
$ws.Range("J9").Value = 0.06880588804507781
$ws.Range("O9").Value = 0.6643867693241158
$ws.Range("P9").Value = 0.6643867693241158
$ws.Range("Q9").Value = 203.6529438677684
$ws.Range("R9").Value = 1832.876494809916
$ws.Range("S9").Value = 0.04571372166874604
$ws.Range("T9").Value = 0.04571372166874605
$ws.Range("G10").Value = 2.079068666666667
$ws.Range("H10").Value = 6.237206
$ws.Range("I10").Value = 0.0688058880450778
$ws.Range("J10").Value = 0.06880588804507781
$ws.Range("M10").Value = 48.20552666666666
$ws.Range("N10").Value = 144.61658
$ws.Range("O10").Value = 0.326960996476427
$ws.Range("P10").Value = 0.326960996476427
$ws.Range("Q10").Value = 100.2226000528311
$ws.Range("R10").Value = 902.00340047548
$ws.Range("S10").Value = 0.02249684171866411
$ws.Range("T10").Value = 0.02249684171866412
